$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.009.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +8.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.509.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +11.75%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '190.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +13.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '548.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.499.51'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +11.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.605'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.73%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.629'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.150'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +18.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.62'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.074.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +11.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.515.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +12.00%  '
$ws.Range('E17').Value = '  +5.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.083.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +8.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.991'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '423.99'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.46'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '652.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.63'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.07%  '
$ws.Range('E34').Value = '  +7.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '59.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0811'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +20.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.387'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.140'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.04%  '
$ws.Range('E41').Value = '  +16.59%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.998.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +16.33%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0415'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.04%  '
$ws.Range('E48').Value = '  +4.50%  '
$ws.Range('E49').Value = '  +7.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +17.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '140.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.79%  '
